# "Grace 3rd updates 0804" - populate the previously-empty "s0" sheet
# (male/female survival-type rates for ages 6-22) in the datainput_bot
# workbook. The sheet starts completely empty (just <sheetData/>), so this
# writes the header row plus the 17 data rows described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("s0")

# Header row: B1/C1 hold the "male"/"female" column labels. Leading "'" so
# they pick up the same quote-prefixed style as the row labels below (the
# target file stores these header cells with a quotePrefix style too).
$ws.Range("B1").Value = "'male"
$ws.Range("C1").Value = "'female"

# Row label (col A, age as text) + male/female rate pair (col B/C) for
# ages 6..22. The age labels are entered with a leading "'" so they are
# stored as text (e.g. "6"), not numbers - matching the source data.
$data = @(
    @("6",  0.027095304550216,  0.02600645374266),
    @("7",  0.0269967348904,    0.025924883055),
    @("8",  0.02657923282358,   0.02552244412122),
    @("9",  0.02616173075676,   0.02512000518744),
    @("10", 0.02574422868994,   0.02471756625366),
    @("11", 0.02532672662312,   0.02431512731988),
    @("12", 0.0249092245563,    0.0239126883861),
    @("13", 0.01998235557356,   0.020238378757836),
    @("14", 0.01963774221982,   0.019885039038852),
    @("15", 0.01929312886608,   0.019531699319868),
    @("16", 0.01189897721506,   0.0151726030368),
    @("17", 0.0116825725874,    0.01489306488),
    @("18", 0.00453967276,      0.00697572612),
    @("19", 0.00444810352,      0.00684189924),
    @("20", 0.00435653428,      0.00670807236),
    @("21", 0.00426496504,      0.00657424548),
    @("22", 0.0041733958,       0.0064404186)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
